# Auto-generated edit script: updates Odin_Profits market-price/profit
# columns (H:N) across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# per the scheduled-runner price refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2130.0852
$ws.Range("J17").Value = 2130.0852
$ws.Range("L17").Value = 6390.2556
$ws.Range("N17").Value = -6726.2556
$ws.Range("H33").Value = 428
$ws.Range("I33").Value = 456.16666
$ws.Range("K33").Value = 456.16666
$ws.Range("M33").Value = -227.16666
$ws.Range("H51").Value = 7821.3125
$ws.Range("I51").Value = 6376.8
$ws.Range("J51").Value = 8477.909
$ws.Range("K51").Value = 6376.8
$ws.Range("L51").Value = 8477.909
$ws.Range("M51").Value = -5892.8
$ws.Range("N51").Value = -9445.909
$ws.Range("H69").Value = 10374.5
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10748
$ws.Range("H70").Value = 7500
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -30540
$ws.Range("H72").Value = 10374.5
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35736
$ws.Range("H73").Value = 7500
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -31872
$ws.Range("H86").Value = 40003036
$ws.Range("J86").Value = 3709.2
$ws.Range("L86").Value = 3709.2
$ws.Range("N86").Value = -5955.2
$ws.Range("H89").Value = 40003036
$ws.Range("J89").Value = 3709.2
$ws.Range("L89").Value = 18546
$ws.Range("N89").Value = -29778
$ws.Range("H127").Value = 5806.4736
$ws.Range("I127").Value = 3293.6667
$ws.Range("J127").Value = 8068
$ws.Range("K127").Value = 9881.000100000001
$ws.Range("L127").Value = 24204
$ws.Range("M127").Value = -4921.000100000001
$ws.Range("N127").Value = -34124
$ws.Range("H132").Value = 514048.8
$ws.Range("J132").Value = 9944.637000000001
$ws.Range("L132").Value = 29833.911
$ws.Range("N132").Value = -34893.911
$ws.Range("H137").Value = 4454.6875
$ws.Range("I137").Value = 1372.875
$ws.Range("J137").Value = 7536.5
$ws.Range("K137").Value = 4118.625
$ws.Range("L137").Value = 22609.5
$ws.Range("M137").Value = -1568.625
$ws.Range("N137").Value = -27709.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11221.444
$ws.Range("I2").Value = 5972.25
$ws.Range("J2").Value = 15420.8
$ws.Range("K2").Value = 5972.25
$ws.Range("L2").Value = 15420.8
$ws.Range("M2").Value = -5859.25
$ws.Range("N2").Value = -15646.8
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 14821524
$ws.Range("I32").Value = 14930289
$ws.Range("J32").Value = 13910609
$ws.Range("K32").Value = 14930289
$ws.Range("L32").Value = 13910609
$ws.Range("M32").Value = -14930002
$ws.Range("N32").Value = -13911183
$ws.Range("H45").Value = 1805.0667
$ws.Range("I45").Value = 1707.1666
$ws.Range("K45").Value = 1707.1666
$ws.Range("M45").Value = -1330.1666
$ws.Range("H102").Value = 1554.4286
$ws.Range("I102").Value = 1249.3334
$ws.Range("K102").Value = 1249.3334
$ws.Range("M102").Value = 372.6666
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72774
$ws.Range("H110").Value = 5998.1934
$ws.Range("I110").Value = 3871.5386
$ws.Range("K110").Value = 3871.5386
$ws.Range("M110").Value = -1826.5386
$ws.Range("H116").Value = 11221.444
$ws.Range("I116").Value = 5972.25
$ws.Range("J116").Value = 15420.8
$ws.Range("K116").Value = 5972.25
$ws.Range("L116").Value = 15420.8
$ws.Range("M116").Value = -3678.25
$ws.Range("N116").Value = -20008.8
$ws.Range("H122").Value = 3403.7021
$ws.Range("I122").Value = 2624.9697
$ws.Range("K122").Value = 7874.909100000001
$ws.Range("M122").Value = -5424.909100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11221.444
$ws.Range("I3").Value = 5972.25
$ws.Range("J3").Value = 15420.8
$ws.Range("K3").Value = 5972.25
$ws.Range("L3").Value = 15420.8
$ws.Range("M3").Value = -5858.25
$ws.Range("N3").Value = -15648.8
$ws.Range("H20").Value = 10991166
$ws.Range("I20").Value = 14287866
$ws.Range("J20").Value = 2166
$ws.Range("K20").Value = 14287866
$ws.Range("L20").Value = 2166
$ws.Range("M20").Value = -14287619
$ws.Range("N20").Value = -2660
$ws.Range("H105").Value = 3038.0625
$ws.Range("I105").Value = 3233
$ws.Range("K105").Value = 3233
$ws.Range("M105").Value = -1486
$ws.Range("H131").Value = 59778.5
$ws.Range("J131").Value = 59778.5
$ws.Range("L131").Value = 59778.5
$ws.Range("N131").Value = -69858.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8063.0835
$ws.Range("I31").Value = 16361.077
$ws.Range("J31").Value = 4980.971
$ws.Range("K31").Value = 16361.077
$ws.Range("L31").Value = 4980.971
$ws.Range("M31").Value = -16066.077
$ws.Range("N31").Value = -5570.971
$ws.Range("H34").Value = 8063.0835
$ws.Range("I34").Value = 16361.077
$ws.Range("J34").Value = 4980.971
$ws.Range("K34").Value = 16361.077
$ws.Range("L34").Value = 4980.971
$ws.Range("M34").Value = -16159.077
$ws.Range("N34").Value = -5384.971
$ws.Range("H68").Value = 61666.625
$ws.Range("J68").Value = 61666.625
$ws.Range("L68").Value = 61666.625
$ws.Range("N68").Value = -63164.625
$ws.Range("H71").Value = 61666.625
$ws.Range("J71").Value = 61666.625
$ws.Range("L71").Value = 184999.875
$ws.Range("N71").Value = -192487.875
$ws.Range("H109").Value = 24752.666
$ws.Range("I109").Value = 12129.5
$ws.Range("J109").Value = 49999
$ws.Range("K109").Value = 12129.5
$ws.Range("L109").Value = 49999
$ws.Range("M109").Value = -11089.5
$ws.Range("N109").Value = -52079
$ws.Range("H122").Value = 2869.7026
$ws.Range("I122").Value = 2445.5833
$ws.Range("J122").Value = 3652.6924
$ws.Range("K122").Value = 7336.749899999999
$ws.Range("L122").Value = 10958.0772
$ws.Range("M122").Value = -4886.749899999999
$ws.Range("N122").Value = -15858.0772
$ws.Range("H131").Value = 74666.664
$ws.Range("I131").Value = 35000
$ws.Range("J131").Value = 94500
$ws.Range("K131").Value = 35000
$ws.Range("L131").Value = 94500
$ws.Range("M131").Value = -29960
$ws.Range("N131").Value = -104580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 104.14286
$ws.Range("I7").Value = 97.25
$ws.Range("J7").Value = 113.333336
$ws.Range("K7").Value = 291.75
$ws.Range("L7").Value = 340.000008
$ws.Range("M7").Value = -179.75
$ws.Range("N7").Value = -564.000008
$ws.Range("H56").Value = 6975.5386
$ws.Range("I56").Value = 6975.5386
$ws.Range("K56").Value = 6975.5386
$ws.Range("M56").Value = -6445.5386
$ws.Range("H113").Value = 1490.25
$ws.Range("J113").Value = 1677.8334
$ws.Range("L113").Value = 5033.5002
$ws.Range("N113").Value = -9373.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5308.1304
$ws.Range("I102").Value = 4009.9656
$ws.Range("K102").Value = 4009.9656
$ws.Range("M102").Value = -2387.9656
$ws.Range("H122").Value = 3298.6758
$ws.Range("I122").Value = 2535.353
$ws.Range("J122").Value = 11949.667
$ws.Range("K122").Value = 7606.059
$ws.Range("L122").Value = 35849.001
$ws.Range("M122").Value = -5156.059
$ws.Range("N122").Value = -40749.001
$ws.Range("H126").Value = 33340242
$ws.Range("I126").Value = 71431110
$ws.Range("J126").Value = 10730.5
$ws.Range("K126").Value = 214293330
$ws.Range("L126").Value = 32191.5
$ws.Range("M126").Value = -214290860
$ws.Range("N126").Value = -37131.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H40").Value = 5657.273
$ws.Range("I40").Value = 4916.778
$ws.Range("K40").Value = 4916.778
$ws.Range("M40").Value = -4780.778
$ws.Range("H43").Value = 490834.25
$ws.Range("J43").Value = 699375
$ws.Range("L43").Value = 699375
$ws.Range("N43").Value = -699761
$ws.Range("H46").Value = 33334560
$ws.Range("I46").Value = 1218
$ws.Range("K46").Value = 1218
$ws.Range("M46").Value = -1030
$ws.Range("H122").Value = 6663333.5
$ws.Range("I122").Value = 9987998
$ws.Range("K122").Value = 29963994
$ws.Range("M122").Value = -29961544
$ws.Range("H132").Value = 7538.207
$ws.Range("I132").Value = 6912.1113
$ws.Range("J132").Value = 8562.727999999999
$ws.Range("K132").Value = 20736.3339
$ws.Range("L132").Value = 25688.184
$ws.Range("M132").Value = -18206.3339
$ws.Range("N132").Value = -30748.184
$ws.Range("H136").Value = 40009764
$ws.Range("I136").Value = 5961.8667
$ws.Range("K136").Value = 17885.6001
$ws.Range("M136").Value = -15335.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29851
$ws.Range("H107").Value = 10527308
$ws.Range("I107").Value = 13334355
$ws.Range("J107").Value = 882.25
$ws.Range("K107").Value = 40003065
$ws.Range("L107").Value = 2646.75
$ws.Range("M107").Value = -40001145
$ws.Range("N107").Value = -6486.75
$ws.Range("H126").Value = 3526.1
$ws.Range("I126").Value = 2016
$ws.Range("K126").Value = 6048
$ws.Range("M126").Value = -3578
$ws.Range("H132").Value = 7314.0566
$ws.Range("I132").Value = 6185.3125
$ws.Range("K132").Value = 18555.9375
$ws.Range("M132").Value = -16025.9375
